# "KB - persona" content refresh (Milestone2 persona slide).
#
# The whole deck is a single slide whose body text lives in a handful of
# shapes; the edit only touches run-level text in two of them:
#
#   Shape 3 "Content Placeholder 5" (Personal Information column)
#     - "Internet " / "usage"                     -> merge to "Internet usage"
#     - "Reddit, Google, " / "WordPress"           -> merge to "Reddit, Google, WordPress"
#     - "University " / "IT, " / "network administration, " / "coding" /
#       ", " / "database design."                  -> merge to 3 runs:
#           "University IT, " | "network administration, coding" | ", database design."
#
#   Shape 6 "TextBox 10" (Business Objectives column)
#     - "We want Chris to"                         -> "We want " / "Patrick to"
#     - "Feel e" / "ncouraged " / "about finding job " -> merge to
#           "Feel encouraged about finding job "
#
# We use the DrawingML TextRange2.Characters(start, length) sub-range API so
# that only the affected characters are rewritten; the host engine then
# re-splits/re-merges the underlying <a:r> runs around the edit, leaving
# every other run (and its rPr) untouched -- exactly mirroring what PowerPoint
# itself does when you select a span of text across run boundaries and
# retype it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Set-Span($shape, [string]$find, [string]$replace) {
    $tr2 = $shape.TextFrame2.TextRange
    $full = $tr2.Text
    $idx = $full.IndexOf($find)
    if ($idx -ge 0) {
        $sub = $tr2.Characters($idx + 1, $find.Length)
        $sub.Text = $replace
    }
}

# --- Shape 3: "Content Placeholder 5" ---------------------------------
$infoShape = $s.Shapes.Item(3)

# "Internet " + "usage" -> single run "Internet usage"
Set-Span $infoShape "Internet usage" "Internet usage"

# "Reddit, Google, " + "WordPress" -> single run "Reddit, Google, WordPress"
Set-Span $infoShape "Reddit, Google, WordPress" "Reddit, Google, WordPress"

# "University " + "IT, " -> single run "University IT, "
Set-Span $infoShape "University IT, " "University IT, "

# "network administration, " + "coding" -> single run "network administration, coding"
Set-Span $infoShape "network administration, coding" "network administration, coding"

# ", " + "database design." -> single run ", database design."
Set-Span $infoShape ", database design." ", database design."

# --- Shape 6: "TextBox 10" ---------------------------------------------
$objShape = $s.Shapes.Item(6)

# "We want Chris to" -> "We want " (unchanged run) + "Patrick to" (new run)
Set-Span $objShape "Chris to" "Patrick to"

# "Feel e" + "ncouraged " + "about finding job " -> single run
Set-Span $objShape "Feel encouraged about finding job " "Feel encouraged about finding job "
